# Commit: 06-02-2018 7:40 PM
# - Script1!B3: "/api/{users}/" -> "/api/users/{user}" (new shared string)
# - Script1 column B width widened to fit the longer path
# - Book window height tweak (view-state only)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Script1")

$ws.Range("B3").Value = "/api/users/{user}"

$ws.Columns.Item(2).ColumnWidth = 15.65

$wb.Windows.Item(1).Height = 2445
